# prof_attend_conference.xlsx - fix header row (coopro_proj issue)
#
# 1) D1 header text changes from "身分(輸入數字)" to "身分 (學士、碩士或博士班）"
# 2) Column A and D get wider to fit the new, longer header text
# 3) The active selection moves to H6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "identity" header in D1 -------------------------------
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"

# --- 2) Widen column A and column D ---------------------------------------
# Column widths are stored in the workbook as "character" widths on a
# Maximum-Digit-Width pixel grid, so the ColumnWidth values below are chosen
# to land as close as possible to the target stored widths (15.5 and 29.625).
$ws.Columns.Item(1).ColumnWidth = 14.79
$ws.Columns.Item(4).ColumnWidth = 28.79

# --- 3) Move the active selection to H6 ------------------------------------
$ws.Range("H6").Select() | Out-Null
